$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet2")

# Season1 Resign Contract Fix - update resigned salary expectations
$ws.Range("B4").Value = 72
$ws.Range("C5").Value = 71

$ws.Range("B83").Value = 75
$ws.Range("C84").Value = 74

$ws.Range("B87").Value = 75
$ws.Range("C88").Value = 74

# Reset view: scroll back to top and select B2 (was scrolled to A79 / C96 selected)
$ws.Range("B2").Select()
